$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix B72: it currently holds the text "3" - convert it to an actual number 3
$ws.Cells.Item(72, 2).Value = 3

# Add new row 73 with the new annotation entry
$ws.Cells.Item(73, 1).Value = "Ruilin"
$ws.Cells.Item(73, 2).NumberFormat = "@"
$ws.Cells.Item(73, 2).Value = "3"
$ws.Cells.Item(73, 2).Style = "Normal"
$ws.Cells.Item(73, 3).Value = "无"
$ws.Cells.Item(73, 4).Value = "DIS"
$ws.Cells.Item(73, 5).Value = "WRI"
$ws.Cells.Item(73, 6).Value = "a447d1c6-38e7-4648-9ebb-727dbeed5375"
$ws.Cells.Item(73, 7).Value = "SkFAWax0-_annotated.xlsx"
$ws.Cells.Item(73, 8).Value = "We can of course remove this part without taking away nothing from the paper's clarity, technical novelty and experimental success."
